$d = $word.ActiveDocument

# 1) "Югорский Государственный Университет" -> "ЮГУ" (exact case, 4 occurrences;
#    leaves the lowercase "Югорский государственный университет" inside
#    "ФГБОУ ВО «Югорский государственный университет»" untouched)
$r = $d.Content
$r.Find.Execute("Югорский Государственный Университет", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "ЮГУ", 2) | Out-Null

# 2) "Ханты-Мансийск, 1111 г." -> "Ханты-Мансийск, 432 г."
$r = $d.Content
$r.Find.Execute("1111", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "432", 2) | Out-Null

# 3) All "02.05.2024" dates -> "04.05.2024" (5 occurrences across the doc)
$r = $d.Content
$r.Find.Execute("02.05.2024", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "04.05.2024", 2) | Out-Null

# 4) "Кукмберу" -> "Gfdgd"
$r = $d.Content
$r.Find.Execute("Кукмберу", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Gfdgd", 2) | Out-Null

# 5) "Чехова 16" -> "Чехова 17"
$r = $d.Content
$r.Find.Execute("Чехова 16", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Чехова 17", 2) | Out-Null

# 6) Expand the list of qualities demonstrated during the internship
$r = $d.Content
$r.Find.Execute("Внимательность, Любознательность, Пунктуальность", $true, $false, $false, $false, $false, `
                 $true, 1, $false, `
                 "Внимательность, Инициативность, Командная работа, Любознательность, Находчивость, Пунктуальность, Своевременность, Стрессоустойчивость, Упорство", `
                 2) | Out-Null

# 7) Final grade "5" -> "4" (use surrounding context so only the grade changes)
$r = $d.Content
$r.Find.Execute("оценивается на «5»", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "оценивается на «4»", 2) | Out-Null

# 8) Remove the two "не посещал занятия" list paragraphs entirely (paragraph + mark)
$toRemove = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("не посещал занятия")) {
        $toRemove += $i
    }
}
for ($j = $toRemove.Count - 1; $j -ge 0; $j--) {
    $idx = $toRemove[$j]
    $d.Paragraphs.Item($idx).Range.Delete()
}
